$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 17:28:55"
$wsZhCn.Range("H2").Value = "2016-03-20 17:29:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 17:29:04"
$wsDeDe.Range("H2").Value = "2016-03-20 17:29:50"
